$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I6").Value = 19.36574074074073
$ws.Range("N6").Value = 2.499849249524808
$ws.Range("O6").Value = 2.870901733221348

$ws.Range("I7").Value = 1.925925925925943
$ws.Range("N7").Value = 2.004846509671994
$ws.Range("O7").Value = 2.229613377609108

$ws.Range("I8").Value = 1.925925925925943
